# Apply the "test cases fixed and finished" edits to Sheet1.
# Package durations were shortened from 2-3 months down to 1 month each,
# labels/test descriptions were renamed to match, the coupon-input question
# was clarified, and the resulting formulas were left in place so they
# recalculate against the new inputs (only the G9 formula actually needed
# to drop its "-20" coupon deduction to match the new expected output).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: Green package w/ extra data, now 1 month (was 3) ---
$ws.Range("A12").Value = "Green Package 1 months with additional 5 GB data"
$ws.Range("B12").Value = "green"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = "yes"

# --- Row 11: Purple package, now 1 month (was 2) ---
$ws.Range("A11").Value = "Purple Package 1 months"
$ws.Range("B11").Value = "purple"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = "no"

# --- Row 10: Blue package, now 1 month (was 2) ---
$ws.Range("A10").Value = "Blue Package 1 months"
$ws.Range("B10").Value = "blue"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = "no"

# --- Row 9: Green package, now 1 month (was 2) ---
$ws.Range("A9").Value = "Green Package 1 months"
$ws.Range("B9").Value = "green"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = "yes"

# Row 9's "cost after coupon" formula no longer subtracts the $20 coupon
# (one month at $49.99 doesn't clear the $75 coupon threshold).
$ws.Range("G9").Formula = "=C2*C9"

# --- Header row 8: clarify the coupon input question ---
$ws.Range("E8").Value = "input: does user have coupon (is it vaild)"

# --- Last selected cell when the sheet was saved ---
[void]$ws.Range("G16").Select()

$wb.Save()
